$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New applicant rows (12:16, 22.10.2024 batch).
$rows = @(
    @{ r = 60; A = "Normurodova Fotima Asaddin qizi"; B = "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik"; C = "AD3672634"; D = "656"; E = "Samarqand viloyati"; F = "Paxtachi tumani"; G = "998934112173"; H = "21-10-2024" },
    @{ r = 61; A = "Husainova Laylo Allayorovna"; B = "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik"; C = "AD4912996"; D = "657"; E = "Buxoro viloyati"; F = "Peshku tumani"; G = "998907158671"; H = "21-10-2024" },
    @{ r = 62; A = "Mamatova Nigora Yo'ldosh qizi"; B = "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik"; C = "AC0870202"; D = "658"; E = "Surxondaryo viloyati"; F = "Sherobod tumani"; G = "998937092606"; H = "22-10-2024" }
)

# Columns D (Shartnoma raqam) and G (Telefon raqam) look numeric; force text
# formatting before assigning so the digits aren't auto-converted to a
# number, then strip the formatting back off so no residual style is left
# on the cell (every other data cell in the sheet uses the default style).
$numericLookingCols = @("D", "G")

foreach ($row in $rows) {
    $r = $row.r
    foreach ($col in @("A", "B", "C", "D", "E", "F", "G", "H")) {
        $cell = $ws.Range("$col$r")
        if ($numericLookingCols -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value = $row[$col]
            $cell.ClearFormats()
        } else {
            $cell.Value = $row[$col]
        }
    }
}
